# Edit script: update commonness-index confidence interval values (columns H/I),
# apply number-format styling to select cells, set column widths, selection, and page margins.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.12436080724000931
$ws.Range("I2").Value = 0.52366214990615845
$ws.Range("H3").Value = 0.12436080724000931
$ws.Range("I3").Value = 0.52366214990615845
$ws.Range("H4").Value = 0.10100854188203812
$ws.Range("I4").Value = 0.50216466188430786
$ws.Range("H5").Value = 0.083859793841838837
$ws.Range("I5").Value = 0.49338066577911377
$ws.Range("H6").Value = 0.76485836505889893
$ws.Range("I6").Value = 0.9034649133682251
$ws.Range("H7").Value = 0.75770026445388794
$ws.Range("I7").Value = 0.89753812551498413
$ws.Range("H8").Value = 0.75443828105926514
$ws.Range("I8").Value = 0.89349114894866943
$ws.Range("H9").Value = 0.76360994577407837
$ws.Range("I9").Value = 0.89504629373550415
$ws.Range("H10").Value = 0.22846108675003052
$ws.Range("I10").Value = 1.1147490739822388
$ws.Range("H11").Value = 0.083164937794208527
$ws.Range("I11").Value = 0.95710194110870361
$ws.Range("H12").Value = 0.12617737054824829
$ws.Range("I12").Value = 0.92468816041946411
$ws.Range("H13").Value = 0.14696630835533142
$ws.Range("I13").Value = 0.951973557472229
$ws.Range("H14").Value = 0.060346517711877823
$ws.Range("I14").Value = 0.47596254944801331
$ws.Range("H15").Value = 0.090785764157772064
$ws.Range("I15").Value = 0.513588547706604
$ws.Range("H16").Value = 0.06515975296497345
$ws.Range("I16").Value = 0.49540045857429504
$ws.Range("H17").Value = 0.060741964727640152
$ws.Range("I17").Value = 0.49687516689300537
$ws.Range("H18").Value = -0.45975089073181152
$ws.Range("I18").Value = 0.26661381125450134
$ws.Range("H19").Value = -0.32624819874763489
$ws.Range("I19").Value = 0.57752424478530884
$ws.Range("H20").Value = -0.23454101383686066
$ws.Range("I20").Value = 0.57059204578399658
$ws.Range("H21").Value = -0.23456914722919464
$ws.Range("I21").Value = 0.62982988357543945
$ws.Range("H22").Value = -0.018348162993788719
$ws.Range("I22").Value = 0.14830964803695679
$ws.Range("H23").Value = -0.011208095587790012
$ws.Range("I23").Value = 0.1521250307559967
$ws.Range("H24").Value = -0.0095542902126908302
$ws.Range("I24").Value = 0.14343246817588806
$ws.Range("H25").Value = -0.0097298407927155495
$ws.Range("I25").Value = 0.13876868784427643
$ws.Range("H26").Value = 0.21055229008197784
$ws.Range("I26").Value = 0.30660146474838257
$ws.Range("H27").Value = 0.20738489925861359
$ws.Range("I27").Value = 0.30938127636909485
$ws.Range("H28").Value = 0.20548927783966064
$ws.Range("I28").Value = 0.30796360969543457
$ws.Range("H29").Value = 0.20986177027225494
$ws.Range("I29").Value = 0.30794167518615723
$ws.Range("H30").Value = 4.1820697784423828
$ws.Range("I30").Value = 6.0359458923339844
$ws.Range("H31").Value = 4.1502089500427246
$ws.Range("I31").Value = 6.0185770988464355
$ws.Range("H32").Value = 4.1483917236328125
$ws.Range("I32").Value = 6.0523786544799805
$ws.Range("H33").Value = 4.1599140167236328
$ws.Range("I33").Value = 5.9788436889648438
$ws.Range("H34").Value = -0.0078555736690759659
$ws.Range("I34").Value = 0.042446907609701157
$ws.Range("H35").Value = -0.0068858270533382893
$ws.Range("I35").Value = 0.044343017041683197
$ws.Range("H36").Value = -0.0086450716480612755
$ws.Range("I36").Value = 0.041958168148994446
$ws.Range("H37").Value = -0.0082412064075469971
$ws.Range("I37").Value = 0.042179439216852188
$ws.Range("H38").Value = -0.18132717907428741
$ws.Range("I38").Value = 0.80455845594406128
$ws.Range("H39").Value = -0.14665378630161285
$ws.Range("I39").Value = 0.82931327819824219
$ws.Range("H40").Value = -0.14474603533744812
$ws.Range("I40").Value = 0.78618878126144409
$ws.Range("H41").Value = -0.14344799518585205
$ws.Range("I41").Value = 0.73258233070373535
$ws.Range("H42").Value = 0.020239254459738731
$ws.Range("I42").Value = 0.11567499488592148
$ws.Range("H43").Value = 0.029953978955745697
$ws.Range("I43").Value = 0.1232246533036232
$ws.Range("H44").Value = 0.026755779981613159
$ws.Range("I44").Value = 0.11859242618083954
$ws.Range("H45").Value = 0.026878176257014275
$ws.Range("I45").Value = 0.1158953458070755
$ws.Range("H46").Value = 0.20910288393497467
$ws.Range("I46").Value = 0.38475227355957031
$ws.Range("H47").Value = 0.20744460821151733
$ws.Range("I47").Value = 0.39385616779327393
$ws.Range("H48").Value = 0.20269899070262909
$ws.Range("I48").Value = 0.39619565010070801
$ws.Range("H49").Value = 0.20428825914859772
$ws.Range("I49").Value = 0.39683151245117188
$ws.Range("H50").Value = 4.0995664596557617
$ws.Range("I50").Value = 6.6781644821166992
$ws.Range("H51").Value = 4.0646967887878418
$ws.Range("I51").Value = 6.7185921669006348
$ws.Range("H52").Value = 4.040229320526123
$ws.Range("I52").Value = 6.8003106117248535
$ws.Range("H53").Value = 4.0494141578674316
$ws.Range("I53").Value = 6.7551913261413574
$ws.Range("H54").Value = 0.0014932226622477174
$ws.Range("I54").Value = 0.043444715440273285
$ws.Range("H55").Value = 0.0025658800732344389
$ws.Range("I55").Value = 0.046721633523702621
$ws.Range("H56").Value = 0.00086158595513552427
$ws.Range("I56").Value = 0.044726207852363586
$ws.Range("H57").Value = 0.0015181305352598429
$ws.Range("I57").Value = 0.044215723872184753
$ws.Range("H58").Value = 0.028615837916731834
$ws.Range("I58").Value = 0.55835682153701782
$ws.Range("H59").Value = 0.070418104529380798
$ws.Range("I59").Value = 0.6069299578666687
$ws.Range("H60").Value = 0.057976618409156799
$ws.Range("I60").Value = 0.58492070436477661
$ws.Range("H61").Value = 0.053307969123125076
$ws.Range("I61").Value = 0.55913043022155762

# Number formats for specific cells/ranges (per diff: new numFmt styles 0.00 and 0.000)
$ws.Range("D26").NumberFormat = "0.00"
$ws.Range("E26").NumberFormat = "0.00"
$ws.Range("D30:D33").NumberFormat = "0.000"

# Column widths
$ws.Columns("A:B").ColumnWidth = 12.66667
$ws.Columns("C:C").ColumnWidth = 5.16667

# Selection matches the authored file's last-saved selection
$ws.Range("D30:D33").Select() | Out-Null

# Page margins (values are specified in points; 72 points per inch)
$ws.PageSetup.LeftMargin = 0.7 * 72
$ws.PageSetup.RightMargin = 0.7 * 72
$ws.PageSetup.TopMargin = 0.75 * 72
$ws.PageSetup.BottomMargin = 0.75 * 72
$ws.PageSetup.HeaderMargin = 0.3 * 72
$ws.PageSetup.FooterMargin = 0.3 * 72
